$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 3.797283712024409
$ws.Range("C2").Value = 0.2129906763635745
$ws.Range("D2").Value = 0.03761025079620595
$ws.Range("E2").Value = 0.04996238046304025
$ws.Range("F2").Value = 6.637023374826981
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1855533544095849
$ws.Range("L2").Value = 0.3688527315500565
$ws.Range("B3").Value = 3.743424325023682
$ws.Range("C3").Value = 0.1987409035744179
$ws.Range("D3").Value = 0.03293470210834926
$ws.Range("E3").Value = 0.05007737893529629
$ws.Range("F3").Value = 6.472368299668346
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1839621642882392
$ws.Range("L3").Value = 0.3694447346956977
$ws.Range("B4").Value = 3.713426736063468
$ws.Range("C4").Value = 0.1902298304181329
$ws.Range("D4").Value = 0.03005861569111801
$ws.Range("E4").Value = 0.05015316951044058
$ws.Range("F4").Value = 6.373076997328241
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1830008451609295
$ws.Range("L4").Value = 0.3700680739760784
$ws.Range("B5").Value = 3.701973054973735
$ws.Range("C5").Value = 0.1868208218203335
$ws.Range("D5").Value = 0.02888500651508252
$ws.Range("E5").Value = 0.0501853602850193
$ws.Range("F5").Value = 6.333061975047315
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1826129447341991
$ws.Range("L5").Value = 0.3703874293019069
$ws.Range("B6").Value = 3.700117665120672
$ws.Range("C6").Value = 0.1862583239030187
$ws.Range("D6").Value = 0.02869002677719834
$ws.Range("E6").Value = 0.05019078446745795
$ws.Range("F6").Value = 6.326444289193944
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.182548763553779
$ws.Range("L6").Value = 0.3704444044263653
$ws.Range("B7").Value = 3.713269150185738
$ws.Range("C7").Value = 0.190183615936661
$ws.Range("D7").Value = 0.03004279472020244
$ws.Range("E7").Value = 0.05015359835705824
$ws.Range("F7").Value = 6.372535540803767
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1829955983605025
$ws.Range("L7").Value = 0.3700721163522687
$ws.Range("B8").Value = 3.778074096969078
$ws.Range("C8").Value = 0.2080275093989883
$ws.Range("D8").Value = 0.03599906398513042
$ws.Range("E8").Value = 0.05000095820831885
$ws.Range("F8").Value = 6.579870528619068
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1850013987654933
$ws.Range("L8").Value = 0.3690029195562801
$ws.Range("B9").Value = 3.929640425085211
$ws.Range("C9").Value = 0.2449407732085263
$ws.Range("D9").Value = 0.0476492681283105
$ws.Range("E9").Value = 0.04974262354301551
$ws.Range("F9").Value = 7.001160671531068
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.1890637728075859
$ws.Range("L9").Value = 0.3689688983776023
$ws.Range("B10").Value = 4.056096055097328
$ws.Range("C10").Value = 0.2732782103547038
$ws.Range("D10").Value = 0.05620750908670402
$ws.Range("E10").Value = 0.04957765890121091
$ws.Range("F10").Value = 7.320188397898704
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1921337321437875
$ws.Range("L10").Value = 0.3702038565431991
$ws.Range("B11").Value = 4.116943542335889
$ws.Range("C11").Value = 0.2864443681191915
$ws.Range("D11").Value = 0.0601043852699803
$ws.Range("E11").Value = 0.0495079728600174
$ws.Range("F11").Value = 7.467506660872232
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.1935503407979482
$ws.Range("L11").Value = 0.3710399777545348
$ws.Range("B12").Value = 4.140465708502916
$ws.Range("C12").Value = 0.2914704472292158
$ws.Range("D12").Value = 0.06158085589467532
$ws.Range("E12").Value = 0.04948235255408107
$ws.Range("F12").Value = 7.523616571598325
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.1940897751049206
$ws.Range("L12").Value = 0.3713960968963619
$ws.Range("B13").Value = 4.135378370721071
$ws.Range("C13").Value = 0.2903861848172653
$ws.Range("D13").Value = 0.06126283113196962
$ws.Range("E13").Value = 0.04948783620794028
$ws.Range("F13").Value = 7.511517763687607
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.1939734632978443
$ws.Range("L13").Value = 0.3713176427162637
$ws.Range("B14").Value = 4.118869077745273
$ws.Range("C14").Value = 0.2868570538063295
$ws.Range("D14").Value = 0.06022583759821032
$ws.Range("E14").Value = 0.04950584967610983
$ws.Range("F14").Value = 7.472116320027339
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.1935946597188263
$ws.Range("L14").Value = 0.3710684839818015
$ws.Range("B15").Value = 4.108819325644731
$ws.Range("C15").Value = 0.2847006353734969
$ws.Range("D15").Value = 0.05959076266032071
$ws.Range("E15").Value = 0.04951698343870492
$ws.Range("F15").Value = 7.448024209400785
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.1933630249478497
$ws.Range("L15").Value = 0.3709210124719959
$ws.Range("B16").Value = 4.052186617698339
$ws.Range("C16").Value = 0.2724233747548794
$ws.Range("D16").Value = 0.05595293697051318
$ws.Range("E16").Value = 0.04958232066947033
$ws.Range("F16").Value = 7.310605624656546
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.1920415665983199
$ws.Range("L16").Value = 0.370154738410335
$ws.Range("B17").Value = 4.01829714440521
$ws.Range("C17").Value = 0.2649627157584575
$ws.Range("D17").Value = 0.05372239029539116
$ws.Range("E17").Value = 0.04962377355306979
$ws.Range("F17").Value = 7.226871022434295
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1912361122478643
$ws.Range("L17").Value = 0.3697549504016138
$ws.Range("B18").Value = 3.999117340065482
$ws.Range("C18").Value = 0.2606974098871717
$ws.Range("D18").Value = 0.05243976932828787
$ws.Range("E18").Value = 0.04964812051981482
$ws.Range("F18").Value = 7.178914792432494
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1907747185415687
$ws.Range("L18").Value = 0.3695508210357872
$ws.Range("B19").Value = 3.992676980278873
$ws.Range("C19").Value = 0.2592576708901788
$ws.Range("D19").Value = 0.05200554426004089
$ws.Range("E19").Value = 0.04965645067408897
$ws.Range("F19").Value = 7.162712705277841
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.19061881854158
$ws.Range("L19").Value = 0.3694861392223032
$ws.Range("B20").Value = 4.02187237158779
$ws.Range("C20").Value = 0.2657542330462377
$ws.Range("D20").Value = 0.05395979959142494
$ws.Range("E20").Value = 0.04961930863836495
$ws.Range("F20").Value = 7.235763362737998
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.1913216586355233
$ws.Range("L20").Value = 0.3697948361024004
$ws.Range("B21").Value = 4.123705192434045
$ws.Range("C21").Value = 0.2878925442060449
$ws.Range("D21").Value = 0.06053040357393513
$ws.Range("E21").Value = 0.0495005378489155
$ws.Range("F21").Value = 7.483680625061538
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.1937058413397921
$ws.Range("L21").Value = 0.3711405955454978
$ws.Range("B22").Value = 4.19306119409049
$ws.Range("C22").Value = 0.3025967111972818
$ws.Range("D22").Value = 0.06482951492927214
$ws.Range("E22").Value = 0.04942739151961795
$ws.Range("F22").Value = 7.647599656653426
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.1952815536459056
$ws.Range("L22").Value = 0.37225038042817
$ws.Range("B23").Value = 4.155787212222151
$ws.Range("C23").Value = 0.2947270187350455
$ws.Range("D23").Value = 0.06253446271823293
$ws.Range("E23").Value = 0.04946602208923689
$ws.Range("F23").Value = 7.559937092225027
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.1944389275098359
$ws.Range("L23").Value = 0.3716369810959463
$ws.Range("B24").Value = 4.020255065615629
$ws.Range("C24").Value = 0.2653963136114612
$ws.Range("D24").Value = 0.05385246764842577
$ws.Range("E24").Value = 0.04962132562153987
$ws.Range("F24").Value = 7.231742563504042
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1912829779117331
$ws.Range("L24").Value = 0.3697767236867122
$ws.Range("B25").Value = 3.885998975846519
$ws.Range("C25").Value = 0.234744377173115
$ws.Range("D25").Value = 0.0444990311174962
$ws.Range("E25").Value = 0.04980813781521609
$ws.Range("F25").Value = 6.885560666850381
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1879502952349839
$ws.Range("L25").Value = 0.3687570434240541
